$d = $word.ActiveDocument

# Apply the topic shifts in reverse order so that a later replacement's
# "new" text never gets re-matched by an earlier replacement's "old" text.

$d.Content.Find.Execute("Augmented Reality", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Mobile games for learning", 2)

$d.Content.Find.Execute("Reading screens", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Augmented Reality", 2)

$d.Content.Find.Execute("1:1 Computing", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Reading screens", 2)

$d.Content.Find.Execute("Situated cognition & embodiment", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1:1 Computing", 2)
